$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.325.16'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '3.001.44'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '563.13'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').Value = '138.91'
$ws.Range('E6').Value = '  +3.83%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.519'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '2.986.62'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('D11').Value = '5.20'
$ws.Range('E11').Value = '  +6.53%  '
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  +1.86%  '
$ws.Range('D13').Value = '0.0000231'
$ws.Range('E13').Value = '  +2.91%  '
$ws.Range('D14').Value = '33.76'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('D16').Value = '7.36'
$ws.Range('E16').Value = '  +6.70%  '
$ws.Range('D17').Value = '3.501.35'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = '3.002.08'
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').Value = '59.322.62'
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('D20').Value = '431.08'
$ws.Range('E20').Value = '  +2.37%  '
$ws.Range('D21').Value = '13.65'
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('D22').Value = '0.722'
$ws.Range('E22').Value = '  +5.00%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '7.14'
$ws.Range('E23').Value = '  +1.73%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '13.51'
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('D25').Value = '80.49'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D27').Value = '2.24'
$ws.Range('E27').Value = '  +11.26%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('D30').Value = '7.88'
$ws.Range('E30').Value = '  +3.54%  '
$ws.Range('D31').Value = '25.77'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('D32').Value = '6.12'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').Value = '0.0999'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '6.00'
$ws.Range('E34').Value = '  +5.69%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +5.88%  '
$ws.Range('D36').Value = '0.0₃0758'
$ws.Range('E36').Value = '  +7.85%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Value = '48.97'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('D40').Value = '2.75'
$ws.Range('E40').Value = '  +5.94%  '
$ws.Range('D41').Value = '409.88'
$ws.Range('E41').Value = '  +8.00%  '
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('D43').Value = '2.775.98'
$ws.Range('E43').Value = '  +3.10%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = '0.253'
$ws.Range('E45').Value = '  +4.02%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '34.96'
$ws.Range('E47').Value = '  +20.93%  '
$ws.Range('D48').Value = '123.57'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('D51').Value = '23.55'
$ws.Range('E51').Value = '  -0.33%  '
